$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: switch the sample/quick-fill row from "Create" to "UpdateByTag" ---
$ws.Cells.Item(5,1).Value = "UpdateByTag"        # A5
$ws.Cells.Item(5,3).ClearContents()              # C5 (was "ddd")
$ws.Cells.Item(5,4).Value = "Brilliant Bulk Update " # D5 (new)
$ws.Cells.Item(5,7).ClearContents()              # G5 (was 53)
$ws.Cells.Item(5,9).Value = "Tag2"               # I5 (new)

# --- Add the "Update Operation" label on existing row 24 ---
$ws.Cells.Item(24,2).Value = "Update Operation"
$ws.Cells.Item(24,2).Font.Bold = $true

# --- Insert two new rows at 25/26 for the UpdateByTag instructions ---
$ws.Rows("25:26").Insert()

$ws.Cells.Item(25,1).Value = "Rem"
$ws.Cells.Item(25,2).Value = "UpdateByTag Operation"
$ws.Cells.Item(25,2).Font.Bold = $true

$ws.Cells.Item(26,1).Value = "Rem"

# --- Append the new UpdateByTag example row at the bottom of the sheet ---
$ws.Cells.Item(36,1).Value = "Rem"
$ws.Cells.Item(36,2).Value = "UpdateByTag"
$ws.Cells.Item(36,3).Value = "Bug"
$ws.Cells.Item(36,4).Value = "New Title to Be added to all Bugs with this tag"
$ws.Cells.Item(36,5).Value = "New Desc."
$ws.Cells.Item(36,10).Value = "Tag to search"

# --- View state: scroll so row 8 is at the top, matching the saved view ---
$ws.Application.ActiveWindow.ScrollRow = 8
